$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.496.85'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '3.128.70'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'215.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = "'621.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.49%  '
$ws.Range("E7").Value = '  +25.72%  '
$ws.Range("D8").Value = "'0.363"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.76%  '
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '3.128.15'
$ws.Range("E10").Value = '  +1.29%  '
$ws.Range("D11").Value = "'0.736"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.53%  '
$ws.Range("E12").Value = '  +5.25%  '
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.28%  '
$ws.Range("D14").Value = "'5.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.70%  '
$ws.Range("D15").Value = "'35.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.38%  '
$ws.Range("D16").Value = '90.243.71'
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("D17").Value = '3.708.24'
$ws.Range("D18").Value = '3.149.49'
$ws.Range("E18").Value = '  +2.24%  '
$ws.Range("D19").Value = "'3.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.66%  '
$ws.Range("D20").Value = "'14.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.90%  '
$ws.Range("D21").Value = "'0.0000213"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.08%  '
$ws.Range("D22").Value = "'461.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.03%  '
$ws.Range("D23").Value = "'9.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.72%  '
$ws.Range("D24").Value = "'5.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.11%  '
$ws.Range("D25").Value = "'95.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +13.66%  '
$ws.Range("D26").Value = "'5.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.98%  '
$ws.Range("D27").Value = "'12.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.23%  '
$ws.Range("D28").Value = '3.305.18'
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").Value = "'0.164"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.18%  '
$ws.Range("D31").Value = "'0.218"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +49.72%  '
$ws.Range("D32").Value = "'9.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.23%  '
$ws.Range("D33").Value = "'26.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.92%  '
$ws.Range("D34").Value = "'517.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("D35").Value = "'0.146"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.56%  '
$ws.Range("E36").Value = '  +5.14%  '
$ws.Range("D37").Value = "'7.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.24%  '
$ws.Range("E38").Value = '  +2.68%  '
$ws.Range("D39").Value = "'3.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.15%  '
$ws.Range("D40").Value = "'0.0909"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +26.43%  '
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").Value = "'0.429"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +15.50%  '
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").Value = "'22.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("D43").Value = "'1.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.95%  '
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").Value = "'2.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.69%  '
$ws.Range("D47").Value = "'0.723"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +19.32%  '
$ws.Range("D48").Value = "'4.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +12.97%  '
$ws.Range("D49").Value = "'150.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.34%  '
$ws.Range("D50").Value = "'1.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.36%  '
$ws.Range("D51").Value = "'45.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.25%  '
